# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) held stale strikeout totals from a previous
# regeneration of the save data. This re-derives/rewrites the per-game K
# values for every data row (rows 2-86) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> newly computed K (strikeouts) value.
$kValues = @{
    2 = 2; 3 = 2; 4 = 1; 5 = 0; 6 = 0; 7 = 1; 8 = 2; 9 = 0
    10 = 5; 11 = 0; 12 = 1; 13 = 0; 14 = 0; 15 = 1; 16 = 1; 17 = 2
    18 = 0; 19 = 2; 20 = 2; 21 = 0; 22 = 2; 23 = 1; 24 = 1; 25 = 1
    26 = 1; 27 = 2; 28 = 0; 29 = 3; 30 = 0; 31 = 1; 32 = 1; 33 = 2
    34 = 2; 35 = 2; 36 = 3; 37 = 1; 38 = 1; 39 = 2; 40 = 1; 41 = 1
    42 = 2; 43 = 0; 44 = 2; 45 = 1; 46 = 0; 47 = 1; 48 = 1; 49 = 0
    50 = 2; 51 = 2; 52 = 1; 53 = 3; 54 = 2; 55 = 1; 56 = 2; 57 = 2
    58 = 2; 59 = 0; 60 = 1; 61 = 0; 62 = 2; 63 = 3; 64 = 3; 65 = 1
    66 = 2; 67 = 2; 68 = 3; 69 = 2; 70 = 2; 71 = 1; 72 = 2; 73 = 0
    74 = 1; 75 = 1; 76 = 1; 77 = 0; 78 = 1; 79 = 2; 80 = 1; 81 = 0
    82 = 2; 84 = 3; 86 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
